$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that used to sit at the very
#    top of the document (inside the Heading1 paragraph). "_GoBack"
#    is hidden from the Bookmarks collection/Count (same as real
#    Word), but it can still be retrieved - and deleted - by name.
# ------------------------------------------------------------------
try {
    $d.Bookmarks.Item("_GoBack").Delete()
} catch {
    # already absent - nothing to do
}

# ------------------------------------------------------------------
# 2) Fix up the "grading for junior players ..." sentence: drop the
#    spell-check wrapped "ages" split and merge the wording, then
#    re-anchor "_GoBack" right before "minimum ages for indicative
#    purposes".
# ------------------------------------------------------------------
$find = $d.Content.Find
$old = "The grading for junior players is summarised in the following table where ages are minimum ages for indicative purposes"
$new = "The grading for junior players is summarised in the following table where ages are minimum ages for indicative purposes"
$find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null

$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd([char]13) -eq $new) { $targetPara = $cand; break }
}

$pStart = $targetPara.Range.Start
$pText = $targetPara.Range.Text
$splitOffset = $pText.IndexOf("minimum")
$splitPos = $pStart + $splitOffset
$splitRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $splitRange) | Out-Null

# ------------------------------------------------------------------
# 3) Swap in the refreshed picture (same embedded image relationship,
#    new crop/size/metadata matching a freshly re-inserted screenshot).
# ------------------------------------------------------------------
$shape = $d.InlineShapes.Item(1)
$picPara = $shape.Range.Paragraphs.Item(1)
$picParaIndex = $picPara.Range.Start
$shape.Delete()

$insertPos = $picParaIndex
$insertRange = $d.Range($insertPos, $insertPos)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
<w:body>
<w:p><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:noProof/><w:color w:val="000000"/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="59375944" wp14:editId="4765138D"><wp:extent cx="6188710" cy="887730"/><wp:effectExtent l="0" t="0" r="2540" b="7620"/><wp:docPr id="1" name="Picture 1" descr="A screenshot of a social media post&#10;&#10;Description automatically generated"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="1" name="IJAGrading.png"/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId9"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="6188710" cy="887730"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@
$insertRange.InsertXML($xml)

Write-Output "done"
